$d = $word.ActiveDocument

$replacements = @(
    @("2024-11-14 Thursday", "2024-11-15 Friday"),
    @("889×7=", "797×7="),
    @("172×9=", "714×3="),
    @("723×8=", "513×3="),
    @("540×8=", "196×5="),
    @("376×7=", "900×6="),
    @("778×3=", "899×7="),
    @("166×2=", "200×5="),
    @("630×2=", "712×7="),
    @("831×7=", "500×2="),
    @("422×9=", "951×8="),
    @("963×9=", "367×8="),
    @("659×7=", "844×7="),
    @("495×2=", "529×5="),
    @("432×6=", "148×8="),
    @("279×7=", "303×9="),
    @("609×2=", "275×7="),
    @("113×4=", "840×5="),
    @("201×2=", "315×5="),
    @("396×3=", "673×3="),
    @("325×9=", "287×2="),
    @("434×2=", "185×5="),
    @("510×5=", "636×9="),
    @("227×7=", "200×8="),
    @("928×7=", "325×6="),
    @("792×9=", "429×2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
